$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.161.00'
$ws.Range("E2").Value = '  +0.45%  '
$ws.Range("D3").Value = '1.679.90'
$ws.Range("E3").Value = '  -0.04%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.31'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.77%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.517'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.55'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.23%  '
$ws.Range("E9").Value = '  +2.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0623'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0891'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.14%  '
$ws.Range("D12").Value = '1.919.73'
$ws.Range("E12").Value = '  +0.15%  '
$ws.Range("D13").Value = '1.690.24'
$ws.Range("E13").Value = '  +0.87%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.18'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.79%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.555'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.63'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.17%  '
$ws.Range("D17").Value = '27.134.58'
$ws.Range("E17").Value = '  +0.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '235.46'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.85'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.15%  '
$ws.Range("D20").Value = '0.0₃0738'
$ws.Range("E20").Value = '  -0.21%  '
$ws.Range("E21").Value = '  +0.13%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.53'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.43%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.52'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.63%  '
$ws.Range("E24").Value = '  -1.90%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.78'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.12%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.40'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.17%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.31'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.95%  '
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("E29").Value = '  +0.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0503'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.80%  '
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.36'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.17%  '
$ws.Range("D33").Value = '1.541.39'
$ws.Range("E33").Value = '  +1.35%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.24'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.93%  '
$ws.Range("E35").Value = '  -3.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.604'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.92%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.943'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.27%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.39'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.22%  '
$ws.Range("E39").Value = '  -1.43%  '
$ws.Range("E40").Value = '  +3.66%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.76'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.43%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '69.02'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.60%  '
$ws.Range("E43").Value = '  +0.19%  '
$ws.Range("E44").Value = '  -0.68%  '
$ws.Range("D45").Value = '1.826.26'
$ws.Range("E45").Value = '  +0.29%  '
$ws.Range("E46").Value = '  +1.38%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '89.86'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.43%  '
$ws.Range("E48").Value = '  +3.88%  '
$ws.Range("E49").Value = '  +6.41%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.22'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.33%  '
$ws.Range("E51").Value = '  -0.03%  '

Write-Output "Updated cryptos list"
